$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 'SANLAM      80904421'
$ws.Range("B4").Value = 'WINDEED   WINDEED D67219 2319J'
$ws.Range("B5").Value = 'MOMENTUM  DE000669339   6793QL'
$ws.Range("B6").Value = 'MOMENTUM    099917130   4543AM'
$ws.Range("B7").Value = 'MOMENTUM    200727015   5888ZL'
$ws.Range("B8").Value = 'MOMENTUM    314386836   7952CM'
$ws.Range("B9").Value = 'HOLLARD   HOL6932470    250101'
$ws.Range("B10").Value = 'MOMENTUM    099917123   4548AM'
$ws.Range("B11").Value = 'MOMENTUM    099917117   4556AM'
$ws.Range("B12").Value = 'EFT WAGES TV D579 B JAN 2025'
$ws.Range("B13").Value = 'EFT WAGES ACC D579 B NOV 202'
$ws.Range("B14").Value = 'WINDEED   WINDEED D67219 2319J'
$ws.Range("B16").Value = 'CLAIMS    DEB_8509060020080'
$ws.Range("B17").Value = 'QNOR001 NORLAND CONS'
$ws.Range("B18").Value = 'EFT WAGES RVP001 TV D'
$ws.Range("B19").Value = 'EFT WAGES RVP009 TV D'
$ws.Range("B20").Value = 'EFT WAGES RVP010 TV D'
$ws.Range("B21").Value = 'EFT WAGES RVP012 TV D'
$ws.Range("B22").Value = 'EFT WAGES RVP013 TV D'
$ws.Range("B23").Value = 'EFT WAGES RVP021 TV D'
$ws.Range("B24").Value = 'EFT WAGES RVP022 TV D'
$ws.Range("B25").Value = 'QPHI001'
$ws.Range("B26").Value = 'DCL PREPAID REFUND'
$ws.Range("B27").Value = 'PAYPROP   DEPOSIT REFUND'
$ws.Range("B29").Value = 'EFT WAGES SEEL D583C ACC DEC'
$ws.Range("B30").Value = 'EFT WAGES SEEL D583C FOOD DE'
$ws.Range("B33").Value = 'BS380000 TRF FROM ABSA CHEQU'
$ws.Range("B34").Value = 'RVP WAG ACCOM JAN2025'
$ws.Range("B35").Value = 'RVP WAGES FOOD JAN2025'
$ws.Range("C38").Value = "'9"
$ws.Range("C38").Style = "Normal"
$ws.Range("B42").Value = 'EFT WAGE D182'
$ws.Range("B43").Value = 'EFT WAGE D182 124040'
$ws.Range("B44").Value = 'D182 EFT WAGE'
$ws.Range("B45").Value = 'EFT WAGE D182 124040'
$ws.Range("B46").Value = 'D182 EFT WAGE'
$ws.Range("B47").Value = 'D0568 QUICK STEP'
$ws.Range("B49").Value = 'D524BEOC14JAN25'
$ws.Range("B50").Value = 'DIV9 BUFFALO CITY MM'
$ws.Range("C50").Value = "'9"
$ws.Range("C50").Style = "Normal"
$ws.Range("B51").Value = 'DIV9 BUFFALO CITY MM'
$ws.Range("C51").Value = "'9"
$ws.Range("C51").Style = "Normal"
$ws.Range("B52").Value = 'D563C EFT WAGE RS MIC'
$ws.Range("B53").Value = 'GARN G DE VOS D535B'
$ws.Range("B54").Value = 'GARN IM TSHABALALA D524B'
$ws.Range("B55").Value = 'EFT WAGE D182'
$ws.Range("B56").Value = 'ABSA PLANT'
$ws.Range("B57").Value = 'ABSA INV TRACKER'
